$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add three new columns (Wins, Losses, Ties) holding the team's season
# record, alongside the existing per-player statistics table.

# Copy the formatting of an existing header cell (bold font, thin border,
# centered/top alignment) onto the three new header cells.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Header labels
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Season record repeated for every player row (2-46)
$ws.Range("AD2:AD46").Value = 108
$ws.Range("AE2:AE46").Value = 54
$ws.Range("AF2:AF46").Value = 0
